$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Add new column I data: header "Beach" and "Shipwreck" for row 2
$ws.Range("I1").Value = "Beach"
$ws.Range("I2").Value = "Shipwreck"

# Update the view selection (scroll back to top-left, select A3 instead of G4)
$ws.Range("A3").Select()

$wb.Save()
